$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weather data rows for December 16 (rows 17-32), continuing the log started in rows 2-16
$data = @(
    ,@(17, "December", 16, "15:30:02", 8.13, 0, 59, 1035, 0.45)
    ,@(18, "December", 16, "16:00:02", 8.13, 0, 57, 1035, 0.45)
    ,@(19, "December", 16, "16:30:03", 7.63, 0, 56, 1035, 0.89)
    ,@(20, "December", 16, "17:00:02", 6.63, 0, 57, 1034, 0.89)
    ,@(21, "December", 16, "17:30:02", 5.31, 0, 59, 1035, 0.45)
    ,@(22, "December", 16, "18:00:03", 3.81, 0, 65, 1035, 0.45)
    ,@(23, "December", 16, "18:30:02", 2.69, 0, 70, 1035, 0.45)
    ,@(24, "December", 16, "19:00:02", 1.88, 0, 73, 1035, 1.88)
    ,@(25, "December", 16, "19:30:02", 0, 0, 76, 1036, 0.45)
    ,@(26, "December", 16, "20:00:02", -0.38, 0, 78, 1036, 0.45)
    ,@(27, "December", 16, "20:30:02", -0.69, 0, 78, 1036, 0.89)
    ,@(28, "December", 16, "21:00:02", -0.88, 0, 82, 1036, 0.89)
    ,@(29, "December", 16, "21:30:02", -1.13, 0, 83, 1036, 1.15)
    ,@(30, "December", 16, "22:00:02", -1.19, 0, 83, 1036, 1.15)
    ,@(31, "December", 16, "22:42:25", -1.5, 0, 83, 1036, 0.56)
    ,@(32, "December", 16, "23:00:06", -1.5, 0, 83, 1036, 0.45)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
